# edit.ps1
# Reproduces (semantically) the OOXML diff described for formats.xlsx:
#  - Sheet1: A6/B6 become formulas FALSE()/TRUE() instead of bare boolean literals.
#  - Selection bookkeeping on Sheet1 (A5) and Sheet2 (D6), matching the diff's
#    <selection> changes.
#  - A handful of now-unused legacy cell styles ("Style1","Untitled1".."Untitled8")
#    are removed, mirroring the cellStyles/cellStyleXfs shrinkage in the diff.
#  - A new trailing worksheet "Sheet4" is added, carrying a single date cell (B2)
#    formatted with a custom "D" (day-of-month) number format plus the
#    accompanying font/fill used for that "Date" style in the target file -
#    this is the actual fdo#55198 regression-test payload.
#  - The originally active sheet/tab (Sheet2) is restored as the active sheet
#    at the end, since adding/touching other sheets shifts focus.

function RGB($r, $g, $b) { return $r + $g * 256 + $b * 65536 }

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: turn the literal boolean cells into FALSE()/TRUE() formulas.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A6").Formula = "=FALSE()"
$ws1.Range("B6").Formula = "=TRUE()"

# Selection bookkeeping (matches diff: Sheet1 selection B5 -> A5).
$ws1.Range("A5").Select()

# ---------------------------------------------------------------------------
# Sheet2: selection bookkeeping (matches diff: Sheet2 selection E3 -> D6).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("D6").Select()

# ---------------------------------------------------------------------------
# Drop the unused legacy "Style1"/"UntitledN" cell styles - the target
# styles.xml collapses cellStyles from 15 entries down to the 6 builtins
# plus a single new "Date" style.
# ---------------------------------------------------------------------------
$unusedStyles = @("Style1", "Untitled1", "Untitled2", "Untitled3", "Untitled4", `
                  "Untitled5", "Untitled6", "Untitled7", "Untitled8")
foreach ($styleName in $unusedStyles) {
    try {
        $wb.Styles.Item($styleName).Delete()
    } catch {
        # Ignore if already absent / not supported - non essential cosmetic cleanup.
    }
}

# ---------------------------------------------------------------------------
# Add the new trailing "Sheet4" with the fdo#55198 regression payload: a date
# serial formatted with a bare "D" (day-of-month) custom number format, using
# the same look (small grey Arial font on a pale-yellow fill) as the "Date"
# cell style introduced in the target workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "Sheet4"

$dateCell = $ws4.Range("B2")
$dateCell.Value = 24470
$dateCell.NumberFormat = "D"
$dateCell.Font.Name = "Arial"
$dateCell.Font.Size = 8
$dateCell.Font.Color = (RGB 64 64 64)
$dateCell.Interior.Color = (RGB 255 255 204)
$dateCell.Interior.PatternColor = (RGB 255 255 255)

# ---------------------------------------------------------------------------
# Restore the originally active sheet (Sheet2 had tabSelected/activeTab) since
# adding Sheet4 and touching other sheets moves the active tab around.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Sheet2").Activate()
